$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This mirrors a manual edit made in the "Cliente" sheet: a new client row
# ("KARLOS MV (W8)") was inserted above row 4, pushing the existing client
# rows (GENIVALDO ... HAREF MACIEL) down by one. Because only columns B:G
# hold client data (column A is a static row counter), the row counter itself
# was left untouched. A handful of existing client names also picked up
# their "(CLIENTE NN/20)" suffix at the same time.
# ---------------------------------------------------------------------------

function Copy-CellFormat($src, $dst) {
    $dst.Font.Bold = $src.Font.Bold
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Color = $src.Font.Color
    $dst.Interior.Color = $src.Interior.Color
    $dst.Borders.LineStyle = $src.Borders.LineStyle
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.WrapText = $src.WrapText
}

# Snapshot the current (pre-edit) B:G contents for the data rows (4..10) so
# the values survive being overwritten while we shift everything down.
$snapshot = @{}
for ($r = 4; $r -le 10; $r++) {
    $snapshot[$r] = @(
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2,
        $ws.Cells.Item($r, 7).Value2
    )
}

# Shift rows 4..10 down into rows 5..11, bottom-up so nothing is clobbered
# before it is read. Formatting is carried along explicitly (cell-by-cell)
# so destination cells that were previously blank (row 11, and F8) end up
# with the same look as the rest of the table instead of staying "empty"
# styled.
for ($r = 10; $r -ge 4; $r--) {
    $destRow = $r + 1
    $vals = $snapshot[$r]

    Copy-CellFormat $ws.Cells.Item($r, 2) $ws.Cells.Item($destRow, 2)
    Copy-CellFormat $ws.Cells.Item($r, 3) $ws.Cells.Item($destRow, 3)
    Copy-CellFormat $ws.Cells.Item($r, 4) $ws.Cells.Item($destRow, 4)
    Copy-CellFormat $ws.Cells.Item($r, 5) $ws.Cells.Item($destRow, 5)
    Copy-CellFormat $ws.Cells.Item($r, 6) $ws.Cells.Item($destRow, 6)
    Copy-CellFormat $ws.Cells.Item($r, 7) $ws.Cells.Item($destRow, 7)
    $ws.Cells.Item($destRow, 4).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($destRow, 5).NumberFormat = "0"

    $ws.Cells.Item($destRow, 2).Value = $vals[0]
    $ws.Cells.Item($destRow, 3).Value = $vals[1]
    $ws.Cells.Item($destRow, 4).Value = $vals[2]
    $ws.Cells.Item($destRow, 5).Value = $vals[3]
    $ws.Cells.Item($destRow, 6).Value = $vals[4]
    $ws.Cells.Item($destRow, 7).Value = $vals[5]
}

# New client inserted at row 4 (row 4 already carries the data-row styling,
# so no formatting fix-up is required here).
$ws.Cells.Item(4, 2).Value = "KARLOS MV (W8)"
$ws.Cells.Item(4, 3).Value = "3a6e38f09a011776bd2d5d0f26815e40"
$ws.Cells.Item(4, 4).Value = 44806
$ws.Cells.Item(4, 5).Value = 60
$ws.Cells.Item(4, 6).Value = "jardson@gmail.com"
$ws.Cells.Item(4, 7).Value = "OK"

# Existing clients renamed to include their position in the "/20" batch.
$ws.Cells.Item(8, 2).Value = "SERGIO (CLIENTE 04/20)"
$ws.Cells.Item(9, 2).Value = "EDSON  BARRETO (CLIENTE 05/20)"
$ws.Cells.Item(10, 2).Value = "WALTER FREIRE (CLIENTE 06/20)"
$ws.Cells.Item(11, 2).Value = "HAREF MACIEL (CLIENTE 07/20)"
